$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "want to go" counts in column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 518
$ws1.Range("F5").Value = 212
$ws1.Range("F7").Value = 227
$ws1.Range("F8").Value = 2219
$ws1.Range("F10").Value = 5460

# Sheet "全部类型" (sheet4): same underlying rows, same updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 518
$ws4.Range("F6").Value = 212
$ws4.Range("F8").Value = 227
$ws4.Range("F11").Value = 2219
$ws4.Range("F13").Value = 5460
